# Generate Report for Archive
#
# This script updates the localization-status workbook:
#  - The "Ready for handoff" status text (used on the Overview sheet as well
#    as the per-language zh-cn / de-de sheets) is changed to "In Translation".
#  - The Status column widths that hosted the old, wider text are narrowed
#    to fit the shorter replacement text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Excel snaps ColumnWidth to whole-pixel increments of the Normal-style font,
# so the narrowest achievable width near the target ~13.41 chars is reached
# by requesting 12.5 (Excel stores it as 13.33.. after the standard pixel
# rounding/padding it always applies to column widths).
$newWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: keep $oldStatus on the left of -eq. Some cells hold the
            # text "True"/"False", which Value2 surfaces as a real [bool];
            # PowerShell's -eq coerces its right operand to the type of the
            # left one, so "$cell.Value2 -eq $oldStatus" would wrongly treat
            # any boolean cell as a match (non-empty string -> $true).
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
                $ws.Columns.Item($cell.Column).ColumnWidth = $newWidth
            }
        }
    }
}
